$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 45207 (2023-10-08) to 45208 (2023-10-09) for every data row (2..34).
for ($row = 2; $row -le 34; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
